$d = $word.ActiveDocument

$d.Content.Find.Execute(": ({{qtdChaves}})", $true, $false, $false, $false, $false,
                         $true, 1, $false, ": {{qtdChaves}}", 2)
